$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("A1").Value = "Run"
$ws.Range("B1").Value = "Test Case Number"
$ws.Range("C1").Value = "discerption"
$ws.Range("D1").Value = "Email"
$ws.Range("E1").Value = "password"

# --- Row 2 ---
$ws.Range("A2").Value = "yes"
$ws.Range("B2").Value = "test case number 1"
$ws.Range("C2").Value = "Check response on entering valid  Credentials(Email and password)"
$ws.Range("D2").Value = "wiasm.mtour@gmail.com"
$ws.Range("E2").Value = 123456789

# --- Row 3 ---
$ws.Range("A3").Value = "yes"
$ws.Range("B3").Value = "test case number 2"
$ws.Range("C3").Value = "Check response on entering valid Email and invalid password"
$ws.Range("D3").Value = "wiasm.mtour@gmail.com"
$ws.Range("E3").Value = 1

# --- Row 4 ---
$ws.Range("A4").Value = "no"
$ws.Range("B4").Value = "test case number 3"
$ws.Range("C4").Value = "Check response on entering invalid Email and valid  password"
$ws.Range("D4").Value = "wiasm.mtour@gmail."
$ws.Range("E4").Value = 123456789

# --- Row 5 ---
$ws.Range("A5").Value = "yes"
$ws.Range("B5").Value = "test case number 4"
$ws.Range("C5").Value = "Check response on entering valid Email and blank  password"
$ws.Range("D5").Value = "wiasm.mtour@gmail.com"
$ws.Range("E5").ClearContents()

# --- Row 6 ---
$ws.Range("A6").Value = "yes"
$ws.Range("B6").Value = "test case number 5"
$ws.Range("C6").Value = "Check response on entering blank Email and vaild  password"
$ws.Range("D6").ClearContents()
$ws.Range("E6").Value = 123456789

# --- Hyperlink on D4 (cell text stays "wiasm.mtour@gmail.", links to the full mailto address) ---
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:wiasm.mtour@gmail.com")

# --- Selection / view ---
$ws.Range("C1").Select()
